$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 3")
$ws.Activate()

$ws.Range("A6").Value = "implementation"
$ws.Range("B6").Value = 3

$ws.Range("A7").Value = "Project presentation"
$ws.Range("B7").Value = 3

$ws.Range("A8").Value = "Excel creation"
$ws.Range("B8").Value = 0.5
$ws.Range("B8").HorizontalAlignment = 1

$ws.Range("B5").Select()
